$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.711557113606915
$ws.Range("C2").Value = 0.02626726547452782
$ws.Range("D2").Value = 1.211696474043208
$ws.Range("E2").Value = 0.9450705515896831
$ws.Range("F2").Value = 0.03521556578416447
$ws.Range("G2").Value = 1.248994055272596

$ws.Range("B3").Value = 6.517836968934673
$ws.Range("C3").Value = 0.07618656717764112
$ws.Range("D3").Value = 8.329172101064454
$ws.Range("E3").Value = 4.758482783963441
$ws.Range("F3").Value = 0.05375514678269903
$ws.Range("G3").Value = 6.413843459826229
